$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New product (YL0010C229) spec rows to append, matching the existing
# (product, test_key, lower_limit, upper_limit) layout.
$newRows = @(
    @("YL0010C229", "MDA", -0.5, 0.5),
    @("YL0010C229", "MDB", -0.5, 0.5),
    @("YL0010C229", "MDL", -0.5, 0.5),
    @("YL0010C229", "STR", 95,   105),
    @("YL0010C229", "TDA", -0.5, 0.5),
    @("YL0010C229", "TDB", -0.5, 0.5),
    @("YL0010C229", "TDL", -0.5, 0.5)
)

$startRow = 9
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Tighten up column widths for the lower/upper limit columns now that the
# new rows have been added.
$ws.Columns.Item(3).ColumnWidth = 9.3333333
$ws.Columns.Item(4).ColumnWidth = 10.3333333

# Reflect where the user ended up after entering the new data.
$ws.Range("D16").Select()
